# Replace the real data-lab project reference with the generic
# demonstration placeholder used for published examples.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "[DL-MAA2016-15]"
$newValue = "[DL-MAA20XX-YY]"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
